$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the changed Price/Volume cells so the values are
# written back as literal strings (matching the source file layout),
# not auto-converted to numbers/percentages by Excel's input parser.
$cells = @("D2","E2","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","E19","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","D27","E27","D39","E39","E40","D41","E41","D42","E42","D43","D44","D45","E45","D46","E46","E47","D48","E49","E50","E51")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "308.49"
$ws.Range("E2").Value = "-0.42%"
$ws.Range("E3").Value = "1.10%"
$ws.Range("D4").Value = "5.140"
$ws.Range("E4").Value = "0.75%"
$ws.Range("D5").Value = "0.08123"
$ws.Range("E5").Value = "-0.46%"
$ws.Range("D6").Value = "1.944"
$ws.Range("E6").Value = "-1.85%"
$ws.Range("D7").Value = "8.147"
$ws.Range("E7").Value = "2.78%"
$ws.Range("D8").Value = "0.9282"
$ws.Range("E8").Value = "-0.58%"
$ws.Range("E9").Value = "-0.17%"
$ws.Range("D10").Value = "0.1936"
$ws.Range("E10").Value = "-0.74%"
$ws.Range("D11").Value = "0.09059"
$ws.Range("E11").Value = "-1.16%"
$ws.Range("D12").Value = "0.03497"
$ws.Range("E12").Value = "0.40%"
$ws.Range("D13").Value = "0.09799"
$ws.Range("E13").Value = "-0.40%"
$ws.Range("D14").Value = "0.001402"
$ws.Range("E14").Value = "-0.66%"
$ws.Range("D15").Value = "0.005881"
$ws.Range("E15").Value = "-2.34%"
$ws.Range("D16").Value = "3.908"
$ws.Range("E16").Value = "9.44%"
$ws.Range("D17").Value = "4.238"
$ws.Range("E17").Value = "0.76%"
$ws.Range("D18").Value = "3.377"
$ws.Range("E18").Value = "-1.99%"
$ws.Range("E19").Value = "0.18%"
$ws.Range("E20").Value = "-0.84%"
$ws.Range("D21").Value = "4.759"
$ws.Range("E21").Value = "-1.47%"
$ws.Range("D22").Value = "0.2427"
$ws.Range("E22").Value = "-1.66%"
$ws.Range("D23").Value = "0.04374"
$ws.Range("E23").Value = "-1.90%"
$ws.Range("D24").Value = "0.001230"
$ws.Range("E24").Value = "-0.81%"
$ws.Range("D25").Value = "0.004846"
$ws.Range("E25").Value = "-0.57%"
$ws.Range("D26").Value = "0.0001302"
$ws.Range("D27").Value = "0.0004005"
$ws.Range("E27").Value = "-9.95%"
$ws.Range("D39").Value = "0.02061"
$ws.Range("E39").Value = "-3.95%"
$ws.Range("E40").Value = "-1.31%"
$ws.Range("D41").Value = "0.007439"
$ws.Range("E41").Value = "-0.52%"
$ws.Range("D42").Value = "0.009800"
$ws.Range("E42").Value = "-2.14%"
$ws.Range("D43").Value = "0.1363"
$ws.Range("D44").Value = "0.002133"
$ws.Range("D45").Value = "0.008499"
$ws.Range("E45").Value = "-12.81%"
$ws.Range("D46").Value = "0.00006405"
$ws.Range("E46").Value = "2.74%"
$ws.Range("E47").Value = "-0.14%"
$ws.Range("D48").Value = "0.002597"
$ws.Range("E49").Value = "-18.90%"
$ws.Range("E50").Value = "-0.14%"
$ws.Range("E51").Value = "-0.14%"
